$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation prompt
$excel.DisplayAlerts = $false

# Remove the MÚSICA and DOODLE sheets, keep only "Calculo de média"
[void]$wb.Worksheets.Item("MÚSICA").Delete()
[void]$wb.Worksheets.Item("DOODLE").Delete()

$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("Calculo de média")
$ws.Activate()

# Update values in row 4: F4 10 -> 0, G4 1 -> 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Update selection on the active sheet
[void]$ws.Range("D11").Select()
